$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 1: new title row above the old header ("classic teseract"), plus a run
# of centered blank cells H1:M1 matching the table's decorative right block.
# The old B1/C1/D1 header content is removed entirely (the header moves to
# row 2 below).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "classic teseract"
$ws.Range("B1").Clear()
$ws.Range("C1").Clear()
$ws.Range("D1").Clear()

$ws.Range("A1").Copy()
$ws.Range("H1:M1").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Row 2: table header (previously row 1), now with a 4th column header "val"
# as plain centered text instead of the numeric-formatted cell it used to be.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "rec_type"
$ws.Range("B2").Value = "val_type"
$ws.Range("C2").Value = "ds"

# D2 previously held a numeric (0.000-formatted) cell; restyle it to plain
# centered text formatting (matching A2/B2/C2) before writing the string.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)
$ws.Range("D2").Value = "val"

$ws.Range("A1").Copy()
$ws.Range("H2:L2").PasteSpecial($xlPasteFormats)
$ws.Range("D3").Copy()
$ws.Range("M2").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Data rows 3-12: rec_type / val_type / ds / value
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 3;  A = "straight_recognition";  B = "accuracy"; C = "ds1"; D = 0 },
    @{ Row = 4;  A = "straight_recognition";  B = "wer";      C = "ds1"; D = 0.954 },
    @{ Row = 5;  A = "straight_recognition";  B = "accuracy"; C = "ds2"; D = 0 },
    @{ Row = 6;  A = "straight_recognition";  B = "wer";      C = "ds2"; D = 1.19 },
    @{ Row = 7;  A = "augmented_recognition"; B = "accuracy"; C = "ds1"; D = 0 },
    @{ Row = 8;  A = "augmented_recognition"; B = "wer";      C = "ds1"; D = 1.045 },
    @{ Row = 9;  A = "with_post_recognition"; B = "accuracy"; C = "ds1"; D = 0 },
    @{ Row = 10; A = "with_post_recognition"; B = "wer";      C = "ds1"; D = 1.045 },
    @{ Row = 11; A = "with_post_recognition"; B = "accuracy"; C = "ds2"; D = 0 },
    @{ Row = 12; A = "with_post_recognition"; B = "wer";      C = "ds2"; D = 1.147 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D

    if ($n -le 11) {
        $ws.Range("A1").Copy()
        $ws.Range("H$n`:L$n").PasteSpecial($xlPasteFormats)
        $ws.Range("D3").Copy()
        $ws.Range("M$n").PasteSpecial($xlPasteFormats)
    }
}

# ---------------------------------------------------------------------------
# Selection, as recorded in the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("O21").Select()
